$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column C keeps its text (inline string) representation instead of
# being auto-converted to a date serial number by Excel's COM layer.
$ws.Range("C2:C81").NumberFormat = "@"

$updates = @(
    @{Row=2; C="01/10/2014"; D=887; E=3.26},
    @{Row=3; C="01/10/2014"; D=2174; E=19.12},
    @{Row=4; C="01/10/2014"; D=1757; E=-20.39},
    @{Row=5; C="01/10/2014"; D=1937; E=-4.11},
    @{Row=6; C="01/10/2014"; D=2441; E=18.38},
    @{Row=7; C="01/10/2014"; D=1605; E=5.04},
    @{Row=8; C="01/10/2014"; D=3050; E=10.27},
    @{Row=9; C="01/10/2014"; D=4010; E=7.42},
    @{Row=10; C="01/10/2015"; D=824; E=-7.1},
    @{Row=11; C="01/10/2015"; D=2062; E=-5.15},
    @{Row=12; C="01/10/2015"; D=1613; E=-8.199999999999999},
    @{Row=13; C="01/10/2015"; D=1908; E=-1.5},
    @{Row=14; C="01/10/2015"; D=2157; E=-11.63},
    @{Row=15; C="01/10/2015"; D=1315; E=-18.07},
    @{Row=16; C="01/10/2015"; D=2823; E=-7.44},
    @{Row=17; C="01/10/2015"; D=3615; E=-9.85},
    @{Row=18; C="01/10/2016"; D=869; E=5.46},
    @{Row=19; C="01/10/2016"; D=2341; E=13.53},
    @{Row=20; C="01/10/2016"; D=1678; E=4.03},
    @{Row=21; C="01/10/2016"; D=1651; E=-13.47},
    @{Row=22; C="01/10/2016"; D=1976; E=-8.390000000000001},
    @{Row=23; C="01/10/2016"; D=1652; E=25.63},
    @{Row=24; C="01/10/2016"; D=3103; E=9.92},
    @{Row=25; C="01/10/2016"; D=4505; E=24.62},
    @{Row=26; C="01/10/2017"; D=1041; E=19.79},
    @{Row=27; C="01/10/2017"; D=1897; E=-18.97},
    @{Row=28; C="01/10/2017"; D=1948; E=16.09},
    @{Row=29; C="01/10/2017"; D=1939; E=17.44},
    @{Row=30; C="01/10/2017"; D=1444; E=-26.92},
    @{Row=31; C="01/10/2017"; D=1441; E=-12.77},
    @{Row=32; C="01/10/2017"; D=2615; E=-15.73},
    @{Row=33; C="01/10/2017"; D=3502; E=-22.26},
    @{Row=34; C="01/10/2018"; D=931; E=-10.57},
    @{Row=35; C="01/10/2018"; D=2012; E=6.06},
    @{Row=36; C="01/10/2018"; D=1484; E=-23.82},
    @{Row=37; C="01/10/2018"; D=1674; E=-13.67},
    @{Row=38; C="01/10/2018"; D=1878; E=30.06},
    @{Row=39; C="01/10/2018"; D=1468; E=1.87},
    @{Row=40; C="01/10/2018"; D=2682; E=2.56},
    @{Row=41; C="01/10/2018"; D=3772; E=7.71},
    @{Row=42; C="01/10/2019"; D=781; E=-16.11},
    @{Row=43; C="01/10/2019"; D=1742; E=-13.42},
    @{Row=44; C="01/10/2019"; D=1461; E=-1.55},
    @{Row=45; C="01/10/2019"; D=1663; E=-0.66},
    @{Row=46; C="01/10/2019"; D=1779; E=-5.27},
    @{Row=47; C="01/10/2019"; D=1384; E=-5.72},
    @{Row=48; C="01/10/2019"; D=2336; E=-12.9},
    @{Row=49; C="01/10/2019"; D=3730; E=-1.11},
    @{Row=50; C="01/10/2020"; D=811; E=3.84},
    @{Row=51; C="01/10/2020"; D=1760; E=1.03},
    @{Row=52; C="01/10/2020"; D=1518; E=3.9},
    @{Row=53; C="01/10/2020"; D=2084; E=25.32},
    @{Row=54; C="01/10/2020"; D=1625; E=-8.66},
    @{Row=55; C="01/10/2020"; D=1632; E=17.92},
    @{Row=56; C="01/10/2020"; D=3082; E=31.93},
    @{Row=57; C="01/10/2020"; D=3896; E=4.45},
    @{Row=58; C="01/10/2021"; D=855; E=5.43},
    @{Row=59; C="01/10/2021"; D=1974; E=12.16},
    @{Row=60; C="01/10/2021"; D=1694; E=11.59},
    @{Row=61; C="01/10/2021"; D=1419; E=-31.91},
    @{Row=62; C="01/10/2021"; D=1797; E=10.58},
    @{Row=63; C="01/10/2021"; D=1175; E=-28},
    @{Row=64; C="01/10/2021"; D=2473; E=-19.76},
    @{Row=65; C="01/10/2021"; D=3964; E=1.75},
    @{Row=66; C="01/10/2022"; D=966; E=12.98},
    @{Row=67; C="01/10/2022"; D=2098; E=6.28},
    @{Row=68; C="01/10/2022"; D=1392; E=-17.83},
    @{Row=69; C="01/10/2022"; D=1557; E=9.73},
    @{Row=70; C="01/10/2022"; D=1612; E=-10.29},
    @{Row=71; C="01/10/2022"; D=1124; E=-4.34},
    @{Row=72; C="01/10/2022"; D=2653; E=7.28},
    @{Row=73; C="01/10/2022"; D=3659; E=-7.69},
    @{Row=74; C="01/10/2023"; D=787; E=-18.53},
    @{Row=75; C="01/10/2023"; D=2083; E=-0.71},
    @{Row=76; C="01/10/2023"; D=1541; E=10.7},
    @{Row=77; C="01/10/2023"; D=1702; E=9.31},
    @{Row=78; C="01/10/2023"; D=1899; E=17.8},
    @{Row=79; C="01/10/2023"; D=1113; E=-0.98},
    @{Row=80; C="01/10/2023"; D=2322; E=-12.48},
    @{Row=81; C="01/10/2023"; D=3391; E=-7.32}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
